$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new forecast row (row 54), mirroring the formatting of the
# previous row (row 53) so the date cell keeps its existing date style.
$ws.Range("A53").Copy($ws.Range("A54"))

$ws.Range("A54").Value2 = 45986
$ws.Range("B54").Value2 = 2025
$ws.Range("C54").Value2 = 2.043309689777173
$ws.Range("D54").Value2 = 2026
$ws.Range("E54").Value2 = 1.199077969291551
